$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Inputs updated (Gas SG and CO2 mole fraction) ---
$ws.Range("B6").Value = 1.51
$ws.Range("C6").Value = 0.99

# --- Default hydrocarbon SG assumption changed from 0.7 to 0.65 ---
$ws.Range("L6").Formula = "= IF(G6<=0,0.65*mwAir,(B6*mwAir - (C6 * H6 + D6 * I6 + E6 * J6 + F6 * K6)) / (1 - C6 - D6 - E6 - F6))"

$q = [char]34
$ws.Range("R4").Formula = "=IF(SUM(C6:F6)>=1," + $q + " Warning: Inert mole fractions leaves no hydrocarbons. Hydrocarbon properties have been defaulted to SG=0.65" + $q + "," + $q + $q + ")"

$ws.Range("R3").Formula = "=IF(L6>16.043," + $q + $q + ",CONCATENATE(" + $q + " Error: Gas SG is too low for inert mole fractions specified. Needs to be at least " + $q + ",ROUND((16.043*(1-C6-D6-E6-F6)+(C6*H6+D6*I6+E6*J6+F6*K6))/mwAir,4)," + $q + " for methane non-inert" + $q + "))"

# --- Selection moved to B7 ---
$ws.Range("B7").Select()
